# Applies the crypto price/volume refresh described in the commit
# "Updated cryptos list on Sun Dec  3 14:58:33 UTC 2023 with GitHub Actions".
#
# Strategy:
#  - Columns D (Price) and E (Volume(1h)) in the sheet are stored as TEXT
#    (t="inlineStr") in the source workbook, not numbers, because some
#    values use '.' as a thousands separator (e.g. "39.722.51") and the
#    volume column keeps literal leading/trailing spaces plus a '%'.
#  - When a new D value looks like a normal decimal number (e.g. "227.05"),
#    Excel's COM layer would silently convert a plain assignment into a
#    real number, dropping the original text semantics and the padding
#    elsewhere. To prevent that, such values are written with a leading
#    apostrophe (Excel's "treat as text" quote-prefix) and the cell's
#    style is reset back to "Normal" right after, so no stray numeric
#    value or button style leaks into the saved file.
#  - Row 48/49 also swap which coin (ARBITRUM / FraxShare) they describe,
#    so B/C/D/E are all rewritten for those two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.668.47'
$ws.Range("E2").Value = '  +2.29%  '
$ws.Range("D3").Value = '2.155.04'
$ws.Range("E3").Value = '  +2.80%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''227.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("E6").Value = '  +1.76%  '
$ws.Range("D7").Value = '''63.18'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.17%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").Value = '''0.391'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.92%  '
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = '''15.85'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("D13").Value = '2.476.34'
$ws.Range("E13").Value = '  +2.59%  '
$ws.Range("D14").Value = '''21.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").Value = '2.159.90'
$ws.Range("E17").Value = '  +2.00%  '
$ws.Range("D18").Value = '39.609.55'
$ws.Range("E18").Value = '  +2.30%  '
$ws.Range("D19").Value = '''71.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = '''229.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.99%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("E24").Value = '  +2.48%  '
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("D26").Value = '''172.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.68%  '
$ws.Range("D27").Value = '''9.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("E28").Value = '  +2.64%  '
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("D30").Value = '''19.80'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.61%  '
$ws.Range("E31").Value = '  +5.37%  '
$ws.Range("E32").Value = '  +1.66%  '
$ws.Range("D33").Value = '''4.57'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.34%  '
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("D35").Value = '''6.92'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.39%  '
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("E37").Value = '  +0.78%  '
$ws.Range("D38").Value = '''3.62'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.21%  '
$ws.Range("D39").Value = '''5.06'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +22.05%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("D41").Value = '''102.70'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.32%  '
$ws.Range("E42").Value = '  -0.46%  '
$ws.Range("D43").Value = '''17.59'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.87%  '
$ws.Range("D44").Value = '1.516.52'
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("D46").Value = '''2.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = '''7.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '''1.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("D50").Value = '''50.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.87%  '
$ws.Range("E51").Value = '  +1.13%  '
